# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect the latest scrape (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml) - cell row -> new value
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value  = 616
$wsExhibition.Range("F5").Value  = 1154
$wsExhibition.Range("F6").Value  = 14294
$wsExhibition.Range("F7").Value  = 16335
$wsExhibition.Range("F18").Value = 101
$wsExhibition.Range("F20").Value = 1248
$wsExhibition.Range("F23").Value = 34
$wsExhibition.Range("F24").Value = 6582
$wsExhibition.Range("F26").Value = 16
$wsExhibition.Range("F29").Value = 5712
$wsExhibition.Range("F33").Value = 4759

# Sheet "全部类型" (sheet4.xml) - row numbers differ slightly from "展览"
# because this sheet contains a few extra rows, but the same F-column updates apply.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 616
$wsAll.Range("F5").Value  = 1154
$wsAll.Range("F6").Value  = 14294
$wsAll.Range("F7").Value  = 16335
$wsAll.Range("F18").Value = 101
$wsAll.Range("F20").Value = 1248
$wsAll.Range("F24").Value = 34
$wsAll.Range("F25").Value = 6582
$wsAll.Range("F27").Value = 16
$wsAll.Range("F32").Value = 5712
$wsAll.Range("F36").Value = 4759
